# Edit script: update the "MicroplasticImages" validation rule row (row 18)
# on the image_explorer_rules sheet, and move the active selection to A18.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("image_explorer_rules")

# Row 18: name, description, valid_example, severity, rule
$ws.Range("A18").Value = "MicroplasticImages"
$ws.Range("B18").Value = "Name of image must be correct"
$ws.Range("C18").Value = "ALGALITA_CW_3_above500_30.jpeg"
$ws.Range("D18").Value = "error"
$ws.Range("E18").Value = "check_exists_in_zip(MicroplasticImages)"

# Move the selection/active cell to A18 (from E18)
$ws.Range("A18").Select()
